# Applies two changes described by the commit diff:
#
# 1. On slide 16, the table's table style is switched from the custom
#    "Table_0" style ({5A2AD623-EA1F-4DAE-92E1-EF2FC855AB0D}, defined in
#    ppt/tableStyles.xml) to the built-in PowerPoint table style
#    {3B8D1971-85ED-480C-8C92-CFD7A5F326A9}.
#
# 2. The presentation's theme ("Integral") and the (separately stored)
#    "Office Theme" swap places: the colour scheme that slides/master
#    currently render with ("Integral") is replaced by the classic
#    "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{3B8D1971-85ED-480C-8C92-CFD7A5F326A9}")
    }
}

# --- 2. Theme colour swap (Integral -> Office Theme) ---------------------------
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgbHex = $officeColors[$i - 1]
    # COM RGB longs are stored as 0x00BBGGRR
    $r = ($rgbHex -band 0xFF0000) -shr 16
    $g = ($rgbHex -band 0x00FF00) -shr 8
    $b = ($rgbHex -band 0x0000FF)
    $comRgb = $b * 65536 + $g * 256 + $r
    $tcs.Colors($i).RGB = $comRgb
}
